$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$ws.Range("B13:D14").NumberFormat = "@"
$ws.Range("B16:D16").NumberFormat = "@"
$ws.Range("B20:D20").NumberFormat = "@"

# Row 13: Enterprises density (per 1000 people)
$ws.Range("B13").Value = "37.76"
$ws.Range("C13").Value = "3.18"
$ws.Range("D13").Value = "40.94"

# Row 14: Employment (% of total)
$ws.Range("B14").Value = "37.07"
$ws.Range("C14").Value = "44.77"
$ws.Range("D14").Value = "81.84"

# Row 16: Enterprises (% of total)
$ws.Range("B16").Value = "92.07"
$ws.Range("C16").Value = "7.75"
$ws.Range("D16").Value = "99.82"

# Row 20: Value added to the economy (% of total)
$ws.Range("B20").Value = "26.57"
$ws.Range("C20").Value = "49.71"
$ws.Range("D20").Value = "76.28"
